# Apply weekly update to the Cebollín sheet:
# Insert two new data rows (one "Primera" and one "Segunda" quality record,
# both dated serial 44944 / 2023-01-18) right after the current row 37,
# pushing all the existing rows 38-78 down by two positions (to 40-80).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 38 (inserting twice at the same index pushes
# the previously-inserted blank row down, leaving two blank rows at 38-39
# and shifting the old row 38 content down to row 40).
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(38).Insert()

# --- New row 38: Cebollín, Primera, 2023-01-18 ---
$ws.Cells.Item(38, 1).Value2 = 7
$ws.Cells.Item(38, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(38, 3).Value2 = "Ñuble"
$ws.Cells.Item(38, 4).Value2 = 44944
$ws.Cells.Item(38, 5).Value2 = 16
$ws.Cells.Item(38, 6).Value2 = 100112037
$ws.Cells.Item(38, 7).Value2 = "Cebollín"
$ws.Cells.Item(38, 8).Value2 = "Sin especificar"
$ws.Cells.Item(38, 9).Value2 = "Primera"
$ws.Cells.Item(38, 10).Value2 = 400
$ws.Cells.Item(38, 11).Value2 = 600
$ws.Cells.Item(38, 12).Value2 = 700
$ws.Cells.Item(38, 13).Value2 = 650
$ws.Cells.Item(38, 14).Value2 = "$/paquete 6 unidades"
$ws.Cells.Item(38, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(38, 16).Value2 = 108
$ws.Cells.Item(38, 17).Value2 = 6
$ws.Cells.Item(38, 18).Value2 = "Hortaliza"

# --- New row 39: Cebollín, Segunda, 2023-01-18 ---
$ws.Cells.Item(39, 1).Value2 = 7
$ws.Cells.Item(39, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(39, 3).Value2 = "Ñuble"
$ws.Cells.Item(39, 4).Value2 = 44944
$ws.Cells.Item(39, 5).Value2 = 16
$ws.Cells.Item(39, 6).Value2 = 100112037
$ws.Cells.Item(39, 7).Value2 = "Cebollín"
$ws.Cells.Item(39, 8).Value2 = "Sin especificar"
$ws.Cells.Item(39, 9).Value2 = "Segunda"
$ws.Cells.Item(39, 10).Value2 = 300
$ws.Cells.Item(39, 11).Value2 = 500
$ws.Cells.Item(39, 12).Value2 = 500
$ws.Cells.Item(39, 13).Value2 = 500
$ws.Cells.Item(39, 14).Value2 = "$/paquete 6 unidades"
$ws.Cells.Item(39, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(39, 16).Value2 = 83
$ws.Cells.Item(39, 17).Value2 = 6
$ws.Cells.Item(39, 18).Value2 = "Hortaliza"
